$wb = $excel.ActiveWorkbook

# --- Helper: update a hyperlink's display text (TextToDisplay) for a given
#     cell address (e.g. "$A$2") on a worksheet, leaving the target Address
#     untouched. We must walk the Hyperlinks collection with foreach since
#     this runtime's .Item(n) indexer on Hyperlinks doesn't resolve properties
#     reliably - only the foreach enumerator does. ---
function Set-HyperlinkDisplay {
    param($ws, [string]$cellAddr, [string]$newDisplay)
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $cellAddr) {
            $h.TextToDisplay = $newDisplay
        }
    }
}

$file50840 = "50840c9f-ebd6-46b3-ba37-fc06ee076493.md"
$fileDc7a0 = "dc7a0273-dac9-469d-8c9c-361251acb6f8.md"

# ===================== Sheet "Overview" =====================
$wsOv = $wb.Worksheets.Item("Overview")

# Row 2 now holds the dc7a0273 file (was 50840c9f), status "In Translation"
$wsOv.Range("A2").Value = $fileDc7a0
$wsOv.Range("B2").Value = "In Translation"
$wsOv.Range("C2").Value = "In Translation"

# Row 3 now holds the 50840c9f file (was dc7a0273), status "Ready for handoff"
$wsOv.Range("A3").Value = $file50840
$wsOv.Range("B3").Value = "Ready for handoff"
$wsOv.Range("C3").Value = "Ready for handoff"

Set-HyperlinkDisplay $wsOv '$A$2' $fileDc7a0
Set-HyperlinkDisplay $wsOv '$A$3' $file50840

# ===================== Sheet "zh-cn" =====================
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhDc7aXlf = "dc7a0273-dac9-469d-8c9c-361251acb6f8.17a5d7b1dbb936cb3b2615b26a3f9d02f311d0f7.zh-cn.xlf"
$zh50840Xlf = "50840c9f-ebd6-46b3-ba37-fc06ee076493.3aafc46c3c43d4a6668076903881bc9086c03c65.zh-cn.xlf"

# Row 2 -> dc7a0273, still "In Translation", handoff file/time unchanged
$wsZh.Range("A2").Value = $fileDc7a0
$wsZh.Range("B2").Value = "In Translation"
$wsZh.Range("C2").Value = $zhDc7aXlf
$wsZh.Range("D2").Value = "2016-03-09 04:20:18"
$wsZh.Range("G2").Value = "0001-01-01 00:00:00"
$wsZh.Range("H2").Value = "Include"

# Row 3 -> 50840c9f, now "Ready for handoff" with a fresh handoff datetime
$wsZh.Range("A3").Value = $file50840
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = $zh50840Xlf
$wsZh.Range("D3").Value = "2016-03-09 04:21:37"
$wsZh.Range("G3").Value = "0001-01-01 00:00:00"
$wsZh.Range("H3").Value = "Include"

Set-HyperlinkDisplay $wsZh '$A$2' $fileDc7a0
Set-HyperlinkDisplay $wsZh '$C$2' $zhDc7aXlf
Set-HyperlinkDisplay $wsZh '$A$3' $file50840
Set-HyperlinkDisplay $wsZh '$C$3' $zh50840Xlf

# ===================== Sheet "de-de" =====================
$wsDe = $wb.Worksheets.Item("de-de")

$deDc7aXlf = "dc7a0273-dac9-469d-8c9c-361251acb6f8.17a5d7b1dbb936cb3b2615b26a3f9d02f311d0f7.de-de.xlf"
$de50840Xlf = "50840c9f-ebd6-46b3-ba37-fc06ee076493.3aafc46c3c43d4a6668076903881bc9086c03c65.de-de.xlf"

# Row 2 -> dc7a0273, still "In Translation", handoff file/time unchanged
$wsDe.Range("A2").Value = $fileDc7a0
$wsDe.Range("B2").Value = "In Translation"
$wsDe.Range("C2").Value = $deDc7aXlf
$wsDe.Range("D2").Value = "2016-03-09 04:20:21"
$wsDe.Range("G2").Value = "0001-01-01 00:00:00"
$wsDe.Range("H2").Value = "Include"

# Row 3 -> 50840c9f, now "Ready for handoff" with a fresh handoff datetime
$wsDe.Range("A3").Value = $file50840
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = $de50840Xlf
$wsDe.Range("D3").Value = "2016-03-09 04:21:39"
$wsDe.Range("G3").Value = "0001-01-01 00:00:00"
$wsDe.Range("H3").Value = "Include"

Set-HyperlinkDisplay $wsDe '$A$2' $fileDc7a0
Set-HyperlinkDisplay $wsDe '$C$2' $deDc7aXlf
Set-HyperlinkDisplay $wsDe '$A$3' $file50840
Set-HyperlinkDisplay $wsDe '$C$3' $de50840Xlf
